$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.106.92"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").Value = "2.527.40"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'596.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "

$ws.Range("D6").Value = "'174.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D8").Value = "'0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "2.525.22"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("D11").Value = "'0.166"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.08%  "

$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("E13").Value = "  -2.14%  "

$ws.Range("D14").Value = "'26.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.72%  "

$ws.Range("D15").Value = "2.986.89"
$ws.Range("E15").Value = "  -0.82%  "

$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("D17").Value = "67.929.00"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").Value = "2.537.84"
$ws.Range("E18").Value = "  -0.93%  "

$ws.Range("D19").Value = "'11.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.75%  "

$ws.Range("D20").Value = "'8.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("D21").Value = "'364.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.45%  "

$ws.Range("D22").Value = "'4.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.32%  "

$ws.Range("D23").Value = "'4.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").Value = "'71.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("E26").Value = "  -3.28%  "

$ws.Range("D27").Value = "'10.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").Value = "2.650.86"
$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("D30").Value = "0.0₃0985"
$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("D31").Value = "'8.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.37%  "

$ws.Range("D32").Value = "'533.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.46%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.60%  "

$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").Value = "'157.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  -2.06%  "

$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").Value = "'18.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("E41").Value = "  -0.58%  "

$ws.Range("D42").Value = "'5.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("E43").Value = "  -1.66%  "

$ws.Range("E44").Value = "  -2.30%  "

$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").Value = "'147.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.78%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0280"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.557"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").Value = "'1.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("E51").Value = "  -1.23%  "
